# Extracted Details From CSVs to Dicts
# The only user-visible change is the header in C1: " Status" (leading
# space) -> "Status" (no leading space), plus the active selection moving
# to C1 (reflecting that the user edited that header cell).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = "Status"

$ws.Range("C1").Select()
